# Apply the "Historias de usuario" content edit to the report document.
$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) The "_GoBack" bookmark currently sits at the end of the title
#    paragraph ("INFORME FINAL DE PROYECTO"). In the target revision it
#    has moved to the end of the first new user-story paragraph, so
#    remove it from its old location first (it will be re-created below).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Locate the final (empty) paragraph of the document - a new
#    "Historias de usuario" section with four user stories is inserted
#    right before it, and the paragraph that used to precede it
#    ("Diseño de interfaces" trailing blank line) becomes a plain,
#    empty List-Paragraph-styled line.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$secondToLast = $d.Paragraphs.Item($count - 1)
$secondToLast.Style = "Prrafodelista"

$lastPara = $d.Paragraphs.Item($count)

$newXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:rPr>
      <w:highlight w:val="yellow"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:highlight w:val="yellow"/>
    </w:rPr>
    <w:t>Historias de usuario</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t>Como administrador, quiero generar un listado de hu&#233;spedes por fecha de llegada y salida, para poder tener un control efectivo de las reservas y facilitar la gesti&#243;n del resort.</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t>Como administrador, quiero gestionar y crear ofertas especiales basadas en las fechas establecidas, para atraer m&#225;s hu&#233;spedes durante temporadas bajas y aumentar la ocupaci&#243;n del resort.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t>Como administrador, quiero poder visualizar todos los registros de los hu&#233;spedes, para realizar un seguimiento de sus preferencias y mejorar la calidad del servicio ofrecido.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t>Como administrador, quiero generar un listado de materiales y equipos disponibles junto con su estado (en uso, mantenimiento, etc.), para gestionar eficazmente los recursos del resort.</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
"@

$lastPara.Range.InsertXML($newXml)
